$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 66, shifting existing rows 66-183 down to 67-184.
$ws.Rows.Item(66).Insert()

# Populate the newly inserted row 66 with the new data point
# (weekly price record for Acelga, Region de Los Lagos origin).
$ws.Range("A66").Value = 4
$ws.Range("B66").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C66").Value = "Los Lagos"
$ws.Range("D66").Value = 44725
$ws.Range("E66").Value = 10
$ws.Range("F66").Value = 100112009
$ws.Range("G66").Value = "Acelga"
$ws.Range("H66").Value = "Sin especificar"
$ws.Range("I66").Value = "Primera"
$ws.Range("J66").Value = 40
$ws.Range("K66").Value = 12000
$ws.Range("L66").Value = 12000
$ws.Range("M66").Value = 12000
$ws.Range("N66").Value = '$/docena de atados (12 kilos)'
$ws.Range("O66").Value = "Región de Los Lagos"
$ws.Range("P66").Value = 1000
$ws.Range("Q66").Value = 12
$ws.Range("R66").Value = "Hortaliza"

# Preserve the date style (numFmt) used throughout column D for the new row.
$ws.Range("D66").NumberFormat = $ws.Range("D67").NumberFormat
